$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder | Antidote
$ws.Range("H6").Value = 124.666664
$ws.Range("I6").Value = 124.666664
$ws.Range("K6").Value = 373.999992
$ws.Range("M6").Value = -261.999992

# Row 8: On the Drip | Eye Drops
$ws.Range("H8").Value = 137.25
$ws.Range("I8").Value = 137.25
$ws.Range("K8").Value = 411.75
$ws.Range("M8").Value = -272.75

# Row 31: Hush Little Wailer | Weak Silencing Potion
$ws.Range("H31").Value = 2092.111
$ws.Range("I31").Value = 1266.125
$ws.Range("K31").Value = 3798.375
$ws.Range("M31").Value = -3568.375

# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 1376.6666
$ws.Range("I33").Value = 315.25
$ws.Range("K33").Value = 315.25
$ws.Range("M33").Value = -86.25

# Row 38: Just Give Him a Serum | Hi-Potion of Strength
$ws.Range("H38").Value = 17763.75
$ws.Range("I38").Value = 16381.286
$ws.Range("K38").Value = 49143.858
$ws.Range("M38").Value = -48771.858

# Row 39: Riches' Brew | Hi-Potion of Mind
$ws.Range("H39").Value = 98.77778000000001
$ws.Range("I39").Value = 98.77778000000001
$ws.Range("K39").Value = 296.33334
$ws.Range("M39").Value = -0.3333400000000211

# Row 40: Stuck in the Moment | Horn Glue
$ws.Range("H40").Value = 1750
$ws.Range("J40").Value = 1000
$ws.Range("L40").Value = 1000
$ws.Range("N40").Value = -1350

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 8569.5
$ws.Range("I132").Value = 8569.5
$ws.Range("K132").Value = 25708.5
$ws.Range("M132").Value = -23178.5

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3374.25
$ws.Range("I138").Value = 3748.5
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 11245.5
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -6105.5
$ws.Range("N138").Value = -19280

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate | Bronze Plate
$ws.Range("H4").Value = 123.5
$ws.Range("I4").Value = 99
$ws.Range("K4").Value = 99
$ws.Range("M4").Value = 17

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 1615.625
$ws.Range("I86").Value = 1615.625
$ws.Range("K86").Value = 1615.625
$ws.Range("M86").Value = -492.625

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 1615.625
$ws.Range("I89").Value = 1615.625
$ws.Range("K89").Value = 8078.125
$ws.Range("M89").Value = -2462.125

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 589.4286
$ws.Range("I94").Value = 356.5
$ws.Range("K94").Value = 356.5
$ws.Range("M94").Value = 94.5

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 2000
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6000
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -3465
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 2020.5
$ws.Range("I31").Value = 1624.6
$ws.Range("K31").Value = 1624.6
$ws.Range("M31").Value = -1329.6

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 2020.5
$ws.Range("I34").Value = 1624.6
$ws.Range("K34").Value = 1624.6
$ws.Range("M34").Value = -1422.6

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 9999.5
$ws.Range("I86").Value = 9999
$ws.Range("K86").Value = 9999
$ws.Range("M86").Value = -8876

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 9999.5
$ws.Range("I89").Value = 9999
$ws.Range("K89").Value = 49995
$ws.Range("M89").Value = -44379

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On | Orange Juice
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 500
$ws.Range("K11").Value = 1500
$ws.Range("M11").Value = -1360

# Row 15: Pretty Enough to Eat | Grilled Carp
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = -340

# Row 50: Moving Up in the World | Rolanberry Cheese
$ws.Range("H50").Value = 113.166664
$ws.Range("I50").Value = 200
$ws.Range("K50").Value = 600
$ws.Range("M50").Value = -119

# Row 53: Rolanberry Fields Forever | Rolanberry Cheese
$ws.Range("H53").Value = 113.166664
$ws.Range("I53").Value = 200
$ws.Range("K53").Value = 600
$ws.Range("M53").Value = -119

# Row 55: Pagan Pastries | Pastry Fish
$ws.Range("H55").Value = 113322.11
$ws.Range("J55").Value = 3499.8
$ws.Range("L55").Value = 10499.4
$ws.Range("N55").Value = -10853.4

# Row 81: It Goes Down Smoothly | Frozen Spirits
$ws.Range("H81").Value = 2015
$ws.Range("I81").Value = 2015
$ws.Range("K81").Value = 6045
$ws.Range("M81").Value = -4922

# Row 84: Quenching the Flame (L) | Frozen Spirits
$ws.Range("H84").Value = 2015
$ws.Range("I84").Value = 2015
$ws.Range("K84").Value = 18135
$ws.Range("M84").Value = -12519

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1607.75
$ws.Range("I131").Value = 943
$ws.Range("K131").Value = 2829
$ws.Range("M131").Value = 2211

$ws = $wb.Worksheets.Item("GSM")
# Row 22: Bad to the Bone | Brass Circlet (Sunstone)
$ws.Range("H22").Value = 23504
$ws.Range("J22").Value = 15000
$ws.Range("L22").Value = 15000
$ws.Range("N22").Value = -16058

# Row 95: Chain of Command | Koppranickel Temple Chain
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2305.8333
$ws.Range("I132").Value = 2367
$ws.Range("K132").Value = 7101
$ws.Range("M132").Value = -4571

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("N7").Value = 0

# Row 16: Saddle Sore | Hard Leather
$ws.Range("H16").Value = 3000
$ws.Range("I16").Value = 3000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2830

# Row 74: Overall, We Blend In | Dhalmelskin Vest
$ws.Range("H74").Value = 12142.714
$ws.Range("I74").Value = 14999.667
$ws.Range("K74").Value = 14999.667
$ws.Range("M74").Value = -14001.667

# Row 77: Eviction Notice (L) | Dhalmelskin Vest
$ws.Range("H77").Value = 12142.714
$ws.Range("I77").Value = 14999.667
$ws.Range("K77").Value = 44999.001
$ws.Range("M77").Value = -40007.001

# Row 81: I Need Your Glove Tonight | Dragonskin Gloves of Healing
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("N81").Value = 0

# Row 84: Halonic Drake Handlers (L) | Dragonskin Gloves of Healing
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("N84").Value = 0

# Row 122: Hell on Leather | Gaja Leather
$ws.Range("H122").Value = 7495
$ws.Range("I122").Value = 9990
$ws.Range("K122").Value = 29970
$ws.Range("M122").Value = -27520

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("WVR")
# Row 3: Trew Enough | Hempen Chausses
$ws.Range("H3").Value = 10000000
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999886

# Row 14: Hat in Hand | Straw Hat
$ws.Range("H14").Value = 40000000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 64: Ribbon of Remembrance | Rainbow Ribbon of Healing
$ws.Range("H64").Value = 50000
$ws.Range("I64").Value = 50000
$ws.Range("K64").Value = 50000
$ws.Range("M64").Value = -49752

# Row 67: The Road Was a Ribbon of Moonlight (L) | Rainbow Ribbon of Healing
$ws.Range("H67").Value = 50000
$ws.Range("I67").Value = 50000
$ws.Range("K67").Value = 50000
$ws.Range("M67").Value = -49142

# Row 96: Skills on Display | Ruby Cotton Cloth
$ws.Range("H96").Value = 466.66666
$ws.Range("I96").Value = 400
$ws.Range("J96").Value = 600
$ws.Range("K96").Value = 400
$ws.Range("L96").Value = 600
$ws.Range("M96").Value = 973
$ws.Range("N96").Value = -3346

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 699
$ws.Range("I122").Value = 699
$ws.Range("K122").Value = 2097
$ws.Range("M122").Value = 353

# Row 125: Color Coated | Almasty Serge Coat of Healing
$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 2086.375
$ws.Range("J126").Value = 444
$ws.Range("L126").Value = 1332
$ws.Range("N126").Value = -6272
